$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the admin product-related API endpoint URIs in column B
# (they now live under /Kmarket/admin/product/... instead of /Kmarket/admin/...)
$ws.Range("B11").Value = "/Kmarket/admin/product/list.do"
$ws.Range("B12").Value = "/Kmarket/admin/product/productDelete.do"
$ws.Range("B13").Value = "/Kmarket/admin/product/productDelete.do"
$ws.Range("B14").Value = "/Kmarket/admin/product/productModify.do"
$ws.Range("B15").Value = "/Kmarket/admin/product/register.do"
$ws.Range("B16").Value = "/Kmarket/admin/product/register.do"
$ws.Range("B17").Value = "/Kmarket/admin/product/category.do"

# Column B grew wider to fit the longer URIs
$ws.Columns.Item(2).ColumnWidth = 40.7

# Update the view: scroll back to the top and move the active selection to B17
$ws.Range("B17").Select() | Out-Null
